$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 922  # was 920
$ws.Cells.Item(6, 6).Value = 158  # was 154
$ws.Cells.Item(7, 6).Value = 940  # was 938
$ws.Cells.Item(8, 6).Value = 735  # was 731
$ws.Cells.Item(9, 6).Value = 181  # was 179
$ws.Cells.Item(12, 6).Value = 769  # was 768
$ws.Cells.Item(13, 6).Value = 254  # was 255
$ws.Cells.Item(14, 6).Value = 550  # was 549
$ws.Cells.Item(16, 6).Value = 1289  # was 1287
$ws.Cells.Item(19, 6).Value = 1081  # was 1072
$ws.Cells.Item(20, 6).Value = 2792  # was 2789
$ws.Cells.Item(21, 6).Value = 1277  # was 1274
$ws.Cells.Item(22, 6).Value = 645  # was 642
$ws.Cells.Item(24, 6).Value = 1240  # was 1239
$ws.Cells.Item(26, 6).Value = 969  # was 967
$ws.Cells.Item(27, 6).Value = 318  # was 317
$ws.Cells.Item(28, 6).Value = 748  # was 689
$ws.Cells.Item(29, 6).Value = 10  # was 8
$ws.Cells.Item(31, 6).Value = 1317  # was 1313

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 510  # was 509

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 717  # was 716

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 717  # was 716
$ws.Cells.Item(7, 6).Value = 510  # was 509
$ws.Cells.Item(8, 6).Value = 510  # was 509
$ws.Cells.Item(12, 6).Value = 922  # was 920
$ws.Cells.Item(13, 6).Value = 158  # was 154
$ws.Cells.Item(15, 6).Value = 940  # was 938
$ws.Cells.Item(16, 6).Value = 735  # was 731
$ws.Cells.Item(17, 6).Value = 181  # was 179
$ws.Cells.Item(25, 6).Value = 769  # was 768
$ws.Cells.Item(26, 6).Value = 254  # was 255
$ws.Cells.Item(27, 6).Value = 550  # was 549
$ws.Cells.Item(29, 6).Value = 1289  # was 1287
$ws.Cells.Item(32, 6).Value = 1081  # was 1072
$ws.Cells.Item(33, 6).Value = 2792  # was 2789
$ws.Cells.Item(34, 6).Value = 1277  # was 1274
$ws.Cells.Item(35, 6).Value = 645  # was 642
$ws.Cells.Item(37, 6).Value = 1240  # was 1239
$ws.Cells.Item(41, 6).Value = 969  # was 967
$ws.Cells.Item(42, 6).Value = 318  # was 317
$ws.Cells.Item(43, 6).Value = 748  # was 689
$ws.Cells.Item(44, 6).Value = 10  # was 8
$ws.Cells.Item(46, 6).Value = 1317  # was 1313
